$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block for year 2018 (rows 44-55), months 01-12
$months = @("01","02","03","04","05","06","07","08","09","10","11","12")
$saida  = @(952,478,521,743,653,461,940,576,530,553,605,1082)
$entrada = @(943,464,537,883,628,535,993,655,492,572,655,1127)

$startRow = 44
for ($i = 0; $i -lt $months.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $months[$i]
    $ws.Cells.Item($r, 2).Value = "2018"
    $ws.Cells.Item($r, 3).Value = $saida[$i]
    $ws.Cells.Item($r, 4).Value = $entrada[$i]
}

$ws.Range("D56").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
